$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update values per corrected IFRS data
$ws.Range("D2").Value = 3775
$ws.Range("E2").Value = 381
$ws.Range("F2").Value = 381
$ws.Range("G2").Value = 523
$ws.Range("H2").Value = 417
$ws.Range("I2").Value = 385
$ws.Range("J2").Value = 32
$ws.Range("K2").Value = 5065
$ws.Range("L2").Value = 965
$ws.Range("M2").Value = 4100
$ws.Range("N2").Value = 3996
$ws.Range("O2").Value = 103
$ws.Range("P2").Value = 65
$ws.Range("Q2").Value = 416
$ws.Range("R2").Value = -452
$ws.Range("S2").Value = -32
$ws.Range("T2").Value = 101
$ws.Range("U2").Value = 314
$ws.Range("V2").Value = 57
$ws.Range("W2").Value = 10.09
$ws.Range("X2").Value = 11.04
$ws.Range("Y2").Value = 10.06
$ws.Range("Z2").Value = 8.51
$ws.Range("AA2").Value = 23.54
$ws.Range("AB2").Value = 5979.55
$ws.Range("AC2").Value = 5928
$ws.Range("AD2").Value = 8.15
$ws.Range("AE2").Value = 61484
$ws.Range("AF2").Value = 0.79
$ws.Range("AG2").Value = 500
$ws.Range("AH2").Value = 1.04
$ws.Range("AI2").Value = 8.43
$ws.Range("AJ2").Value = 6500000

# Row 3: update values per corrected IFRS data
$ws.Range("D3").Value = 3353
$ws.Range("E3").Value = 384
$ws.Range("F3").Value = 384
$ws.Range("G3").Value = 541
$ws.Range("H3").Value = 435
$ws.Range("I3").Value = 402
$ws.Range("J3").Value = 32
$ws.Range("K3").Value = 5421
$ws.Range("L3").Value = 942
$ws.Range("M3").Value = 4479
$ws.Range("N3").Value = 4349
$ws.Range("O3").Value = 129
$ws.Range("P3").Value = 65
$ws.Range("Q3").Value = 365
$ws.Range("R3").Value = -278
$ws.Range("S3").Value = -52
$ws.Range("T3").Value = 61
$ws.Range("U3").Value = 304
$ws.Range("V3").Value = 43
$ws.Range("W3").Value = 11.44
$ws.Range("X3").Value = 12.97
$ws.Range("Y3").Value = 9.640000000000001
$ws.Range("Z3").Value = 8.300000000000001
$ws.Range("AA3").Value = 21.04
$ws.Range("AB3").Value = 6526.14
$ws.Range("AC3").Value = 6192
$ws.Range("AD3").Value = 6.45
$ws.Range("AE3").Value = 66914
$ws.Range("AF3").Value = 0.6
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 1.25
$ws.Range("AI3").Value = 8.08
$ws.Range("AJ3").Value = 6500000

# Row 4: update values per corrected IFRS data
$ws.Range("D4").Value = 3323
$ws.Range("E4").Value = 368
$ws.Range("F4").Value = 368
$ws.Range("G4").Value = 501
$ws.Range("H4").Value = 373
$ws.Range("I4").Value = 347
$ws.Range("J4").Value = 26
$ws.Range("K4").Value = 5820
$ws.Range("L4").Value = 1030
$ws.Range("M4").Value = 4791
$ws.Range("N4").Value = 4641
$ws.Range("O4").Value = 149
$ws.Range("P4").Value = 65
$ws.Range("Q4").Value = 312
$ws.Range("R4").Value = -222
$ws.Range("S4").Value = -33
$ws.Range("T4").Value = 243
$ws.Range("U4").Value = 68
$ws.Range("V4").Value = 48
$ws.Range("W4").Value = 11.08
$ws.Range("X4").Value = 11.23
$ws.Range("Y4").Value = 7.71
$ws.Range("Z4").Value = 6.64
$ws.Range("AA4").Value = 21.5
$ws.Range("AB4").Value = 6978.96
$ws.Range("AC4").Value = 5335
$ws.Range("AD4").Value = 6.78
$ws.Range("AE4").Value = 71406
$ws.Range("AF4").Value = 0.51
$ws.Range("AG4").Value = 500
$ws.Range("AH4").Value = 1.38
$ws.Range("AI4").Value = 9.369999999999999
$ws.Range("AJ4").Value = 6500000

# Row 5: update values per corrected IFRS data
$ws.Range("D5").Value = 3000
$ws.Range("E5").Value = 151
$ws.Range("F5").Value = 151
$ws.Range("G5").Value = 271
$ws.Range("H5").Value = 232
$ws.Range("I5").Value = 232
$ws.Range("J5").Value = -12
$ws.Range("K5").Value = 5847
$ws.Range("L5").Value = 998
$ws.Range("M5").Value = 4849
$ws.Range("N5").Value = 4849
$ws.Range("O5").Value = 135
$ws.Range("P5").Value = 65
$ws.Range("Q5").Value = 254
$ws.Range("R5").Value = -487
$ws.Range("S5").Value = -54
$ws.Range("T5").Value = 466
$ws.Range("U5").Value = -212
$ws.Range("V5").Value = 23
$ws.Range("W5").Value = 5.03
$ws.Range("X5").Value = 7.73
$ws.Range("Y5").Value = 4.89
$ws.Range("Z5").Value = 3.98
$ws.Range("AA5").Value = 20.57
$ws.Range("AB5").Value = 7299.82
$ws.Range("AC5").Value = 3567
$ws.Range("AD5").Value = 10.5
$ws.Range("AE5").Value = 74604
$ws.Range("AF5").Value = 0.5
$ws.Range("AG5").Value = 550
$ws.Range("AH5").Value = 1.47
$ws.Range("AI5").Value = 15.42
$ws.Range("AJ5").Value = 6500000

# Row 6: update values per corrected IFRS data
$ws.Range("D6").Value = 3070
$ws.Range("E6").Value = 17
$ws.Range("F6").Value = 17
$ws.Range("G6").Value = 404
$ws.Range("H6").Value = 388
$ws.Range("I6").Value = 386
$ws.Range("K6").Value = 6579
$ws.Range("L6").Value = 1083
$ws.Range("M6").Value = 5495
$ws.Range("N6").Value = 5174
$ws.Range("P6").Value = 65
$ws.Range("Q6").Value = 103
$ws.Range("R6").Value = -190
$ws.Range("S6").Value = -25
$ws.Range("T6").Value = 494
$ws.Range("U6").Value = -391
$ws.Range("V6").Value = 34
$ws.Range("W6").Value = 0.55
$ws.Range("X6").Value = 12.65
$ws.Range("Y6").Value = 7.71
$ws.Range("Z6").Value = 6.25
$ws.Range("AA6").Value = 19.71
$ws.Range("AB6").Value = 7803.76
$ws.Range("AC6").Value = 5942
$ws.Range("AD6").Value = 4.21
$ws.Range("AE6").Value = 79596
$ws.Range("AF6").Value = 0.31
$ws.Range("AG6").Value = 550
$ws.Range("AH6").Value = 2.2
$ws.Range("AI6").Value = 9.26
$ws.Range("AJ6").Value = 6500000

# Row 7: remove stale data cells (kept only A7:C7)
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8: remove stale data cells (kept only A8:C8)
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9: remove stale data cells (kept only A9:C9)
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
